$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated loading_percent values per row (columns C,D,E,F,G,H,I,M,N)
$data = @{
  2 = @("1.926591911704221", "3.570224562475687", "9.125973891635345", "27.86142683426528", "38.53867165461079", "14.38313679654742", "22.08593602313881", "21.94261168104884", "16.89070536417639")
  3 = @("1.941275389174592", "3.572952874075863", "9.189294839070454", "27.03065380287096", "36.83583778496398", "14.21043604631417", "21.56885043109201", "21.01950731517857", "16.64147149552905")
  4 = @("1.950616509116717", "3.576034341209739", "9.231594291963212", "26.51957206326146", "35.76630386326914", "14.10876578876127", "21.25442451133596", "20.4341085654819", "16.48866825580466")
  5 = @("1.954505465274617", "3.577639989168407", "9.249691438172297", "26.31142762704091", "35.3251952395622", "14.06847766538556", "21.12727378277551", "20.19120892321623", "16.42653034119779")
  6 = @("1.955156214634378", "3.577927627605235", "9.252748375994898", "26.27688420682414", "35.25165570024784", "14.06185815748394", "21.10622587175326", "20.15062381382568", "16.41622243822177")
  7 = @("1.950668622683567", "3.576054583758571", "9.231834875518965", "26.51676394118004", "35.76037515506167", "14.1082177627309", "21.25270546437238", "20.43084985505153", "16.48782961561196")
  8 = @("1.931587569318054", "3.570871266217953", "9.147096968512782", "27.57539869811671", "37.95694524892963", "14.32271060684135", "21.907133038987", "21.62837529658", "16.80476618046452")
  9 = @("1.896726150270188", "3.57201793624667", "9.008073211198219", "29.62815470538437", "42.04490165046567", "14.77589288470462", "23.20579826989639", "23.81644118861089", "17.42490155467812")
  10 = @("1.872635605884013", "3.579946539327695", "8.922493608864498", "31.10307207014173", "44.88129989230801", "15.12573355888835", "24.15767434192853", "25.3120296301076", "17.87551612672585")
  11 = @("1.861998795038497", "3.585126646219454", "8.887159232508747", "31.76328039460046", "46.12988735669257", "15.28789297995735", "24.58793481922511", "25.96574682491624", "18.07858354858234")
  12 = @("1.858016612089235", "3.587316773220157", "8.874296301811951", "32.01147824765835", "46.59632802553199", "15.34967600222985", "24.75029160267172", "26.20930474543771", "18.15513673907561")
  13 = @("1.858872220975433", "3.586834893313005", "8.877043546809677", "31.95810878510154", "46.49616045783145", "15.33635406950501", "24.71535337540438", "26.15702996999849", "18.13866599363185")
  14 = @("1.861670264911306", "3.585302237254405", "8.886090619816493", "31.78373745580364", "46.16839136670072", "15.29296868049837", "24.60130447954809", "25.98586529196221", "18.08488890014181")
  15 = @("1.863390088131973", "3.584393264982847", "8.891699606192272", "31.67668701601708", "45.96678332104124", "15.26644125840483", "24.53136647591959", "25.88049794437987", "18.05190210357788")
  16 = @("1.873337179223931", "3.579639870201847", "8.924875184984653", "31.05968662368924", "44.79882914904496", "15.11519187618977", "24.1294849601571", "25.26875742568144", "17.86220069057502")
  17 = @("1.879521489650827", "3.577128523049937", "8.946148630258802", "30.67823015255357", "44.07135976472072", "15.0231371348839", "23.88210851761651", "24.88653050293409", "17.74528412846969")
  18 = @("1.883108893290466", "3.575832037759918", "8.958723030715614", "30.45782986772814", "43.64903166636628", "14.97047689967993", "23.7395749395951", "24.66418471320875", "17.67786004948907")
  19 = @("1.88432875679154", "3.575418409717827", "8.963038630577616", "30.38304342202576", "43.50538000778711", "14.95269810872716", "23.69127832532214", "24.58847836062022", "17.65500317382604")
  20 = @("1.87886002207132", "3.577380527727054", "8.943849002287866", "30.71894198661656", "44.14920738051055", "15.03290720488105", "23.90846927694796", "24.92747908932545", "17.75774889229937")
  21 = @("1.860847174603793", "3.585746196628037", "8.88341923004687", "31.83500561256772", "46.26484074201756", "15.30570224785946", "24.63482032759551", "26.0362499691036", "18.10069440110671")
  22 = @("1.852758083033133", "3.592546528774343", "8.846940973130454", "32.55377013915983", "47.61025195618658", "15.48615947701853", "25.10611581437446", "26.73757786716215", "18.32279431253196")
  23 = @("1.855457931678815", "3.588794437131537", "8.866134044053595", "32.17120712344974", "46.89570096378213", "15.38966623714718", "24.8549448215375", "26.3654464151804", "18.20446332345951")
  24 = @("1.879158971969227", "3.577266137715749", "8.944887592591741", "30.70053955951187", "44.11402520326618", "15.0284893319245", "23.89655254073003", "24.9089743145263", "17.75211421273173")
  25 = @("1.905887144251007", "3.570475961597038", "9.042775216098388", "29.07737843691318", "40.9662344360572", "14.65010847554891", "22.85406205497198", "23.24327879444686", "17.25772232692745")
}

$cols = @("C","D","E","F","G","H","I","M","N")

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $cols[$i] + $row
    $ws.Range($cell).Value = [double]$vals[$i]
  }
}

$wb.Save()
